$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append two new log rows (116 and 117) to the feed logs sheet
$ws.Range("A116").Value = 115
$ws.Range("B116").Value = 1
$ws.Range("C116").Value = "2024-06-17 09:14:03"
$ws.Range("D116").Value = 200
$ws.Range("E116").Value = 8

$ws.Range("A117").Value = 116
$ws.Range("B117").Value = 2
$ws.Range("C117").Value = "2024-06-17 09:14:04"
$ws.Range("D117").Value = 200
$ws.Range("E117").Value = 1
